# Finished concatenated lang sims table
# Fill in the "Verb" column (F) with the finished verb forms that
# correspond to the verb stems already present in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "stop"
$ws.Range("F3").Value = "bounce"
$ws.Range("F4").Value = "roll"
$ws.Range("F5").Value = "shake"
$ws.Range("F6").Value = "balance"

# Update the selected cell to reflect where the author left off editing.
$ws.Range("G2").Select()
